# Refitting NCDEs to individual patients (for manuscript figure)
# - Adds a new "Label" column (H) marking Control (0) vs MDD (1) rows
# - Updates the refit Prediction/Error values (columns D/E) and the
#   Cross Entropy Loss value (F11) for the first (100-iteration) block

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column H (copy the bold/bordered header style from G1,
# then set the text so the style copy doesn't clobber it)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Label"

# --- Block 1 (Iterations = 100), rows 2-11 ---
$ws.Range("D2").Value = 0.5579882147988364
$ws.Range("E2").Value = 0.5579882147988364
$ws.Range("H2").Value = 0

$ws.Range("D3").Value = 0.470622752987122
$ws.Range("E3").Value = 0.470622752987122
$ws.Range("H3").Value = 0

$ws.Range("D4").Value = 0.6530166547886727
$ws.Range("E4").Value = 0.6530166547886727
$ws.Range("H4").Value = 0

$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.362393591240193
$ws.Range("E6").Value = 0.362393591240193
$ws.Range("H6").Value = 0

$ws.Range("D7").Value = 0.658593943020475
$ws.Range("E7").Value = 0.341406056979525
$ws.Range("H7").Value = 1

$ws.Range("D8").Value = 0.6879180585669117
$ws.Range("E8").Value = 0.3120819414330883
$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 0.5176946047968746
$ws.Range("E9").Value = 0.4823053952031254
$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 0.7486056548250788
$ws.Range("E10").Value = 0.2513943451749212
$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 0.7126048134506341
$ws.Range("E11").Value = 0.2873951865493659
$ws.Range("F11").Value = 0.5425586700439453
$ws.Range("H11").Value = 1

# --- Block 2 (Iterations = 200), rows 12-21: only the new Label column ---
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
